$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 202: Mishkin Energizer ----
$ws.Range("B202").Value = "Mishkin Energizer"
$ws.Range("C202").Value = "YES"
$ws.Range("D202").Value = "NO"
$ws.Range("E202").Value = "String + implementation"
$ws.Range("F202").Value = "Its very hard to explain this , the only thing I learnt from this is that I need to try harder in problems(doesn't mean give more time,  it means think more clearly)"

# ---- Row 203: Path With Minimum Effort (topic filled in later, see below) ----
$ws.Range("B203").Value = "Path With Minimum Effort"
$ws.Range("C203").Value = "YES"
$ws.Range("D203").Value = "NO"
$ws.Range("F203").Value = "You use bfs but greedily , for each node in the current queue we choose the one which has currently lowest effort and thus we need a priority queue."

# ---- Row 204: Find Minimum Time to Reach Last Room II ----
$ws.Range("B204").Value = "Find Minimum Time to Reach Last Room II"
$ws.Range("C204").Value = "YES"
$ws.Range("D204").Value = "YES"
$ws.Range("E204").Value = "Shortest path"
$ws.Range("F204").Value = "This pattern is very useful the bfs modified with the priority queue or set, however here set works but pq gives tle"

# ---- Now fill the topic for row 203 (reuses the "Shortest path" string created above) ----
$ws.Range("E203").Value = "Shortest path"

# ---- Hyperlinks: add in the same order as the target relationship ids (rId90..rId93) ----
# NOTE: the cells already contain the intended display text (set above), so TextToDisplay
# is intentionally omitted -- passing it (even when equal to the existing text) keeps the
# <hyperlink> element out of the saved sheet XML under this runtime.
$ws.Hyperlinks().Add($ws.Range("B202"), "https://codeforces.com/problemset/problem/1257/C2")
$ws.Hyperlinks().Add($ws.Range("B203"), "https://leetcode.com/problems/path-with-minimum-effort/description/")
$ws.Hyperlinks().Add($ws.Range("B200"), "https://www.geeksforgeeks.org/problems/chocolates-pickup/1?utm_source=youtube&utm_medium=collab_striver_ytdescription&utm_campaign=chocolates-pickup")
$ws.Hyperlinks().Add($ws.Range("B204"), "https://leetcode.com/problems/find-minimum-time-to-reach-last-room-ii/description/")

# ---- Restore the "Hyperlink" cell style (font/underline) on the touched cells without ----
# ---- leaving Excel's auto-generated extra style applied to the cell (copy format from an ----
# ---- existing hyperlink cell instead of using .Style, which keeps styles.xml minimal). ----
$ws.Range("B2").Copy()
$ws.Range("B200").PasteSpecial(-4122)
$ws.Range("B202").PasteSpecial(-4122)
$ws.Range("B203").PasteSpecial(-4122)
$ws.Range("B204").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- View state: mirror the sheetView/selection change recorded in the edit ----
$ws.Application.ActiveWindow.ScrollRow = 186
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("F205").Select()
